# Apply "work breakdown" updates: refresh Status column for several tasks,
# shift several end dates (and the resulting Duration formulas recalc
# automatically), add a new "Not applicable" status for the draft change-order
# row, nudge column F's width, and leave the selection on H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (F) updates -------------------------------------------
$ws.Range("F3").Value  = "Done"
$ws.Range("F4").Value  = "Done"
$ws.Range("F5").Value  = "Done"
$ws.Range("F6").Value  = "Ongoing"
$ws.Range("F7").Value  = "Done"
$ws.Range("F9").Value  = "Done"
$ws.Range("F10").Value = "Done"

# --- Date updates (Duration column recalculates automatically) -----------
$ws.Range("D14").Value = 44699
$ws.Range("D15").Value = 44699
$ws.Range("D17").Value = 44699
$ws.Range("D18").Value = 44699

$ws.Range("C19").Value = 44695
$ws.Range("D19").Value = 44699

# New status value for row 19 (draft change-order task), previously blank
$ws.Range("F19").Value = "Not applicable"

# --- Column F width (small cosmetic widening) -----------------------------
$ws.Columns("F").ColumnWidth = 12.83

# --- Selection -------------------------------------------------------------
$ws.Range("H19").Select()
